$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("הגר אגמון", 1),
  @("ליהי בראל", 1),
  @("הילה שולויס", 1),
  @("תאיו ורד", 1),
  @("יולי קזמה", 1),
  @("תאיו ורד", 6),
  @("ליהי בראל", 6),
  @("אביב ואסקז", 1),
  @("דן פימה", 1),
  @("גלי זליג", 1),
  @("אורי שטרנברג", 1),
  @("יהלי דוייב", 1),
  @("ליאם דיין", 1),
  @("מעיין סטרוזר", 1),
  @("יולי יערי תליו", 1),
  @("יהלי גודר", 1),
  @("לינוי קוסטיקה", 1),
  @("שלו דיין", 1),
  @("יהלי גודר", 6),
  @("גלי זליג", 6),
  @("רומי הרשקוביץ", 1),
  @("תומר ששון", 1),
  @("תאיו ורד", 1),
  @("הילה שולויס", 1),
  @("קרן רינת פביאן", 1),
  @("רומי הרשקוביץ", 6),
  @("תאיו ורד", 6),
  @("דן פימה", 1),
  @("עדן ורד מרי", 1),
  @("איתי הראל", 1),
  @("יהלי גודר", 1),
  @("ירון גלפנד", 1),
  @("יולי קזמה", 1),
  @("שלו דיין", 1),
  @("מעיין סטרוזר", 1),
  @("אייל קוטלר", 1),
  @("ליהי בראל", 1),
  @("יולי יערי תליו", 1),
  @("ירון גלפנד", 6),
  @("יהלי גודר", 6),
  @("אביב ואסקז", 1),
  @("עדן ורד מרי", 1),
  @("יהלי דוייב", 1),
  @("ליהי בראל", 1),
  @("הילה שולויס", 1),
  @("איתי הראל", 1),
  @("איתי בסטקר", 1),
  @("ירון גלפנד", 1),
  @("תומר ששון", 1),
  @("יהלי דוייב", 6),
  @("איתי בסטקר", 6),
  @("הגר אגמון", 1),
  @("רומי הרשקוביץ", 1),
  @("אן מרש", 1),
  @("מעיין סטרוזר", 1),
  @("איתי הראל", 1),
  @("יולי קזמה", 1),
  @("שלו דיין", 1),
  @("קרן רינת פביאן", 1),
  @("אורי שטרנברג", 1),
  @("תאיו ורד", 1),
  @("הגר אגמון", 6),
  @("מעיין סטרוזר", 6)
)

$row = 357
foreach ($item in $data) {
  $ws.Cells.Item($row, 1).Value = $item[0]
  $ws.Cells.Item($row, 2).Value = $item[1]
  $row = $row + 1
}

$ws.Range("A377").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 356
$win.ScrollColumn = 1

